$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "27.421.15"
$ws.Range("E2").Value = "  +6.31%  "
Set-TextValue "D3" "1.810.45"
$ws.Range("E3").Value = "  +6.08%  "
$ws.Range("E4").Value = "  +0.10%  "
Set-TextValue "D5" "345.11"
$ws.Range("E5").Value = "  +4.29%  "
Set-TextValue "D6" "1.000"
$ws.Range("E6").Value = "  +0.12%  "
Set-TextValue "D7" "0.3832"
$ws.Range("E7").Value = "  +3.92%  "
Set-TextValue "D8" "50.11"
$ws.Range("E8").Value = "  +3.83%  "
Set-TextValue "D9" "0.3512"
$ws.Range("E9").Value = "  +6.03%  "
Set-TextValue "D10" "1.234"
$ws.Range("E10").Value = "  +5.32%  "
$ws.Range("E11").Value = "  +5.26%  "
Set-TextValue "D12" "1.001"
$ws.Range("E12").Value = "  +0.22%  "
Set-TextValue "D13" "22.53"
$ws.Range("E13").Value = "  +12.48%  "
Set-TextValue "D14" "6.608"
$ws.Range("E14").Value = "  +6.59%  "
Set-TextValue "D15" "7.212"
$ws.Range("E15").Value = "  +4.90%  "
Set-TextValue "D16" "1.808.84"
$ws.Range("E16").Value = "  +6.48%  "
Set-TextValue "D17" "0.00001123"
$ws.Range("E17").Value = "  +5.12%  "
Set-TextValue "D18" "0.06761"
$ws.Range("E18").Value = "  +2.07%  "
Set-TextValue "D19" "86.85"
Set-TextValue "D20" "1.000"
Set-TextValue "D21" "17.77"
$ws.Range("E21").Value = "  +9.73%  "
Set-TextValue "D22" "6.535"
$ws.Range("E22").Value = "  +7.78%  "
Set-TextValue "D24" "27.411.88"
$ws.Range("E24").Value = "  +6.42%  "
Set-TextValue "D25" "2.462"
$ws.Range("E25").Value = "  -0.06%  "
Set-TextValue "D26" "2.679"
$ws.Range("E26").Value = "  +7.92%  "
Set-TextValue "D27" "22.20"
$ws.Range("E27").Value = "  +15.92%  "
Set-TextValue "D28" "1.498"
$ws.Range("E28").Value = "  +14.96%  "
Set-TextValue "D29" "154.41"
$ws.Range("E29").Value = "  +3.23%  "
Set-TextValue "D30" "2.012.29"
$ws.Range("E30").Value = "  +6.49%  "
Set-TextValue "D31" "136.87"
$ws.Range("E31").Value = "  +7.06%  "
Set-TextValue "D32" "6.370"
$ws.Range("E32").Value = "  +6.92%  "
Set-TextValue "D33" "4.079"
$ws.Range("E33").Value = "  -0.73%  "
Set-TextValue "D34" "13.88"
$ws.Range("E34").Value = "  +7.78%  "
Set-TextValue "D35" "0.08819"
$ws.Range("E35").Value = "  +3.57%  "
Set-TextValue "D36" "1.723"
$ws.Range("E36").Value = "  +2.02%  "
Set-TextValue "D37" "5.632"
$ws.Range("E37").Value = "  +5.44%  "
Set-TextValue "D38" "0.7094"
$ws.Range("E38").Value = "  +15.94%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D39" "0.2264"
$ws.Range("E39").Value = "  +6.67%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D40" "0.02418"
$ws.Range("E40").Value = "  +7.29%  "
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D41" "0.06520"
$ws.Range("E41").Value = "  +5.00%  "
Set-TextValue "D42" "8.979"
$ws.Range("E42").Value = "  +5.01%  "
Set-TextValue "D43" "1.287"
$ws.Range("E43").Value = "  +0.95%  "
Set-TextValue "D44" "14.94"
$ws.Range("E44").Value = "  +2.97%  "
Set-TextValue "D45" "0.6587"
$ws.Range("E45").Value = "  +12.77%  "
Set-TextValue "D46" "1.0000"
$ws.Range("E46").Value = "  +0.09%  "
Set-TextValue "D47" "3.987"
$ws.Range("E47").Value = "  +3.66%  "
Set-TextValue "D48" "2.181"
$ws.Range("E48").Value = "  +8.86%  "
Set-TextValue "D49" "132.81"
$ws.Range("E49").Value = "  +4.90%  "
Set-TextValue "D50" "0.07361"
$ws.Range("E50").Value = "  +1.99%  "
$ws.Range("E51").Value = "  +5.30%  "
